$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.587.65'
$ws.Range("E2").Value = '  -2.57%  '

$ws.Range("D3").Value = '2.010.27'
$ws.Range("E3").Value = '  -4.73%  '

$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +0.74%  '

$ws.Range("D5").Value = '332.16'
$ws.Range("E5").Value = '  -3.93%  '

$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("D7").Value = '0.5042'
$ws.Range("E7").Value = '  -3.95%  '

$ws.Range("D8").Value = '0.4259'
$ws.Range("E8").Value = '  -4.16%  '

$ws.Range("D9").Value = '53.61'
$ws.Range("E9").Value = '  -2.22%  '

$ws.Range("D10").Value = '0.09209'
$ws.Range("E10").Value = '  -3.33%  '

$ws.Range("D11").Value = '1.126'
$ws.Range("E11").Value = '  -4.17%  '

$ws.Range("D12").Value = '23.61'
$ws.Range("E12").Value = '  -6.08%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '2.022.81'
$ws.Range("E13").Value = '  -3.88%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '8.123'
$ws.Range("E14").Value = '  -7.66%  '

$ws.Range("D15").Value = '6.549'
$ws.Range("E15").Value = '  -5.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.10'
$ws.Range("E16").Value = '  -5.74%  '

$ws.Range("E17").Value = '  +0.59%  '

$ws.Range("D18").Value = '0.00001126'
$ws.Range("E18").Value = '  -3.79%  '

$ws.Range("D19").Value = '0.06675'
$ws.Range("E19").Value = '  -0.97%  '

$ws.Range("D20").Value = '19.97'
$ws.Range("E20").Value = '  -6.14%  '

$ws.Range("D21").Value = '1.011'
$ws.Range("E21").Value = '  +0.42%  '

$ws.Range("D22").Value = '6.008'
$ws.Range("E22").Value = '  -5.04%  '

$ws.Range("D23").Value = '29.629.46'
$ws.Range("E23").Value = '  -2.54%  '

$ws.Range("E24").Value = '  -5.12%  '

$ws.Range("D25").Value = '2.286'
$ws.Range("E25").Value = '  -1.37%  '

$ws.Range("D26").Value = '159.64'
$ws.Range("E26").Value = '  -2.65%  '

$ws.Range("D27").Value = '20.83'
$ws.Range("E27").Value = '  -5.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.460'
$ws.Range("E28").Value = '  -6.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.350'
$ws.Range("E29").Value = '  -7.77%  '

$ws.Range("D30").Value = '128.68'
$ws.Range("E30").Value = '  -3.68%  '

$ws.Range("D31").Value = '1.064'
$ws.Range("E31").Value = '  -7.49%  '

$ws.Range("D32").Value = '1.595'
$ws.Range("E32").Value = '  -8.61%  '

$ws.Range("D33").Value = '0.09974'
$ws.Range("E33").Value = '  -5.66%  '

$ws.Range("D34").Value = '5.889'
$ws.Range("E34").Value = '  -6.14%  '

$ws.Range("D35").Value = '3.809'

$ws.Range("D36").Value = '9.618'
$ws.Range("E36").Value = '  -8.67%  '

$ws.Range("D37").Value = '0.02475'
$ws.Range("E37").Value = '  -6.08%  '

$ws.Range("D38").Value = '1.328'
$ws.Range("E38").Value = '  -1.53%  '

$ws.Range("D39").Value = '0.06412'
$ws.Range("E39").Value = '  -6.00%  '

$ws.Range("D40").Value = '0.6594'
$ws.Range("E40").Value = '  -6.40%  '

$ws.Range("D41").Value = '11.84'
$ws.Range("E41").Value = '  -6.08%  '

$ws.Range("D42").Value = '0.2081'
$ws.Range("E42").Value = '  -6.90%  '

$ws.Range("E43").Value = '  +0.52%  '

$ws.Range("D44").Value = '0.6378'
$ws.Range("E44").Value = '  -6.98%  '

$ws.Range("D45").Value = '13.65'
$ws.Range("E45").Value = '  -5.65%  '

$ws.Range("D46").Value = '2.218'
$ws.Range("E46").Value = '  -6.16%  '

$ws.Range("D47").Value = '1.288'
$ws.Range("E47").Value = '  -5.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.540'
$ws.Range("E48").Value = '  -3.22%  '

$ws.Range("D49").Value = '0.07017'
$ws.Range("E49").Value = '  -3.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.140'
$ws.Range("E50").Value = '  -5.31%  '

$ws.Range("D51").Value = '0.00000000322'
$ws.Range("E51").Value = '  -6.40%  '
